# Apply cryptos list update (values scraped on Thu Jul 20 13:38:52 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.272.77'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '1.920.12'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'0.8100"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').Value = "'244.47"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = "'0.3253"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.00%  '
$ws.Range('D9').Value = "'27.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.34%  '
$ws.Range('D10').Value = "'0.07257"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.04%  '
$ws.Range('D11').Value = "'0.7905"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.03%  '
$ws.Range('D12').Value = "'0.08090"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').Value = '1.914.88'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = "'5.409"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D15').Value = "'93.99"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.18%  '
$ws.Range('D16').Value = '30.277.89'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = "'14.24"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').Value = "'6.069"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.45%  '
$ws.Range('D19').Value = "'250.17"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').Value = "'0.000007845"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '2.165.62'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = "'8.198"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +20.16%  '
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = "'0.1665"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +17.56%  '
$ws.Range('D26').Value = "'9.497"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('D27').Value = "'167.99"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  +0.65%  '
$ws.Range('E29').Value = '  +6.37%  '
$ws.Range('D30').Value = "'1.389"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('D31').Value = "'1.550"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.40%  '
$ws.Range('D32').Value = "'4.341"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('D33').Value = "'0.05756"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.13%  '
$ws.Range('D34').Value = "'4.153"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('E35').Value = '  +2.49%  '
$ws.Range('D36').Value = "'0.7496"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.76%  '
$ws.Range('D37').Value = "'1.003"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').Value = "'2.728"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('D39').Value = "'0.01960"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('D40').Value = "'2.822"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.55%  '
$ws.Range('D41').Value = "'0.4551"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('D42').Value = "'74.31"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('D43').Value = "'5.979"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.53%  '
$ws.Range('D44').Value = "'0.8525"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.27%  '
$ws.Range('D45').Value = "'1.929"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.93%  '
$ws.Range('D46').Value = "'1.000"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('D47').Value = "'103.69"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('D48').Value = '1.034.93'
$ws.Range('E48').Value = '  +4.96%  '
$ws.Range('B49').Value = 'SynthetixNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D49').Value = "'3.122"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.41%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'9.996"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.33%  '
$ws.Range('D51').Value = "'7.631"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.28%  '
